# Add "hydrogen combined cycle" as a new power-plant type (row 24) on the
# BDPbES sheet, and rename the existing "hydrogen" entry (row 23) to
# "hydrogen combustion turbine".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BDPbES")

# --- Rename row 23's plant type, and give it a (black, vertically
#     centered) look -------------------------------------------------------
$ws.Range("A23").Value2 = "hydrogen combustion turbine"
$ws.Range("A23").Font.Color = 0
$ws.Range("A23").VerticalAlignment = -4108   # xlCenter

# --- Add new row 24: "hydrogen combined cycle" --------------------------
$ws.Range("A24").Value2 = "hydrogen combined cycle"

# Reuse A23's freshly-created format for A24 instead of re-applying the
# font/alignment properties a second time (keeps the styles table tidy).
$ws.Range("A23").Copy()
$ws.Range("A24").PasteSpecial(-4122)   # xlPasteFormats

# Same priority value (2) as every other plant type, with the same
# relative formulas used throughout the rest of the table.
$ws.Range("B24").Value2 = 2
$ws.Range("C24").Formula = '=$B24'
$ws.Range("D24:AK24").Formula = '=$B24'

# --- Match the selection left on the sheet by the author -----------------
$ws.Range("B23:AK24").Select()
